$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

$data = @{
    1 = @(-47.093411805028012, -45.796570636779407, 24.326001492432852, 17.626031549730353, 54.438998573787885, 35.313186192113037, 0.81638124054605699, 0.20978629136840976, 0.059078591904416292)
    2 = @(-38.265634945855453, -37.635360113043362, 38.712456830791211, 23.016927870370008, 105.62959350970402, 35.787042247864107, 0.90498762369666896, 0.034968685182419405, 0.020731388653787819)
    3 = @(-31.571276733407295, -15.939521141709742, 24.791782046323892, 22.854270814486988, 81.88700034240496, 102.24796447119215, 0.0025296428264772169, 0.80934649864970498, 0.61691900283934253)
    4 = @(-93.827219101650854, -75.259887621435894, 20.524946103312175, 24.524406722136476, 22.656434572086848, 34.629856974980676, 0.077412085834651881, 0.6894181029641091, 0.40548889126828125)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = $i + 1
        $ws1.Cells.Item($r, $col).Value = $vals[$i]
    }
}

# Sheet2 holds formulas (=INDEX(Sheet1!...)) that depend on Sheet1 and will
# recalculate automatically, updating the cached <v> values to match.
$wb.Application.CalculateFull()
